# Updates 2014/12-2018/12 (IFRS연결) actuals for rows 2-6 and clears the
# now-unavailable 2019(E)-2021(E) estimate rows 7-9 down to their labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Row 2..6: replace the financial figures with the corrected values
$newValues = @{
    2 = @{ "D" = 2595; "E" = 46; "F" = 46; "G" = 227; "H" = 159; "I" = 159; "J" = 0; "K" = 5104; "L" = 856; "M" = 4248; "N" = 4248; "O" = 0; "P" = 80; "Q" = 137; "R" = -337; "S" = -51; "T" = 2; "U" = 135; "V" = 318; "W" = 1.78; "X" = 6.13; "Y" = 3.8; "Z" = 3.14; "AA" = 20.14; "AB" = 5187.76; "AC" = 994; "AD" = 21.08; "AE" = 26551; "AF" = 0.79; "AG" = 300; "AH" = 1.43; "AI" = 30.18; "AJ" = 16000000 }
    3 = @{ "D" = 2917; "E" = -64; "F" = -64; "G" = 229; "H" = 167; "I" = 167; "J" = 0; "K" = 5101; "L" = 728; "M" = 4373; "N" = 4373; "O" = 0; "P" = 80; "Q" = -246; "R" = 444; "S" = -123; "T" = 1; "U" = -247; "V" = 252; "W" = -2.2; "X" = 5.74; "Y" = 3.89; "Z" = 3.28; "AA" = 16.65; "AB" = 5326.23; "AC" = 1047; "AD" = 15.71; "AE" = 27331; "AF" = 0.6; "AG" = 200; "AH" = 1.22; "AI" = 19.11; "AJ" = 16000000 }
    4 = @{ "D" = 3047; "E" = 59; "F" = 59; "G" = 369; "H" = 277; "I" = 277; "J" = 0; "K" = 5252; "L" = 667; "M" = 4586; "N" = 4586; "O" = 0; "P" = 80; "Q" = 89; "R" = 356; "S" = -168; "T" = 4; "U" = 85; "V" = 121; "W" = 1.93; "X" = 9.1; "Y" = 6.19; "Z" = 5.36; "AA" = 14.54; "AB" = 5626.98; "AC" = 1733; "AD" = 8.34; "AE" = 28660; "AF" = 0.5; "AG" = 250; "AH" = 1.73; "AI" = 14.42; "AJ" = 16000000 }
    5 = @{ "D" = 4131; "E" = 294; "F" = 294; "G" = 336; "H" = 265; "I" = 265; "J" = 0; "K" = 5499; "L" = 703; "M" = 4796; "N" = 4796; "O" = 0; "P" = 80; "Q" = 369; "R" = -95; "S" = -67; "T" = 7; "U" = 361; "V" = 97; "W" = 7.12; "X" = 6.42; "Y" = 5.65; "Z" = 4.93; "AA" = 14.67; "AB" = 5892.77; "AC" = 1657; "AD" = 9.93; "AE" = 29975; "AF" = 0.55; "AG" = 500; "AH" = 3.04; "AI" = 30.17; "AJ" = 16000000 }
    6 = @{ "D" = 3552; "E" = 74; "F" = 74; "G" = 517; "H" = 404; "I" = 404; "K" = 5810; "L" = 712; "M" = 5098; "N" = 5098; "P" = 80; "Q" = -31; "R" = 63; "S" = -106; "T" = 279; "U" = -310; "V" = 74; "W" = 2.09; "X" = 11.37; "Y" = 8.17; "Z" = 7.14; "AA" = 13.97; "AB" = 6296.16; "AC" = 2525; "AD" = 5.11; "AE" = 31862; "AF" = 0.4; "AG" = 500; "AH" = 3.88; "AI" = 19.8; "AJ" = 16000000 }
}

foreach ($row in $newValues.Keys) {
    foreach ($col in $newValues[$row].Keys) {
        $ws.Range("$col$row").Value = $newValues[$row][$col]
    }
}

# Rows 7..9 (2019(E)/2020(E)/2021(E)) no longer have data available;
# clear every figure column but keep the A/B/C labels intact.
$clearCols = @("D", "E", "G", "H", "I", "K", "L", "M", "N", "P", "Q", "R", "S", "T", "U", "W", "X", "Y", "Z", "AA", "AC", "AD", "AE", "AF", "AG", "AH", "AI")
foreach ($row in 7..9) {
    foreach ($col in $clearCols) {
        $ws.Range("$col$row").ClearContents()
    }
}
